# Update the RF column (column I) for rows 26 through 63 with the
# new recalculated Raising Factor value, per "Update of 2025 data and RF changes".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 43.24575

for ($row = 26; $row -le 63; $row++) {
    $ws.Cells.Item($row, 9).Value = $newValue
}
